# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give H1 the same header formatting (bold, centered, bordered) as the
# existing header cells by copying G1's format onto it, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the new column's data values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
